# Add data for 2022-04-23
# Updates the "through" date label and a handful of carjacking counts
# for the current (April 2022) reporting column plus small corrections
# scattered across several neighborhood rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet / update the header label that tracks the "through" date.
$ws.Name = "Through 2022-04-15"
$ws.Range("B1").Value = "April 2022 (through April 15)"

# Austin
$ws.Range("N2").Value = 3
$ws.Range("V2").Value = 1

# North Lawndale
$ws.Range("B4").Value = 2

# New City
$ws.Range("B7").Value = 1

# Woodlawn
$ws.Range("AD18").Value = 1

# Lincoln Park
$ws.Range("N24").Value = 1

# South Shore
$ws.Range("N26").Value = 3
$ws.Range("AD26").Value = 1

# Uptown
$ws.Range("F27").Value = 1

# West Loop
$ws.Range("B28").Value = ""

# Roseland
$ws.Range("V32").Value = 3

# South Deering
$ws.Range("B34").Value = 1
$ws.Range("F34").Value = 2

# Galewood
$ws.Range("R41").Value = 1

# Douglas
$ws.Range("F45").Value = 2

# Lincoln Square
$ws.Range("R74").Value = 1

# North Park
$ws.Range("R80").Value = 1

# South Chicago
$ws.Range("B88").Value = 2
